# Automated map update (2025-07-24 10:12:11)
# Inserts one new record at row 7 (all following rows shift down by one),
# matching the upstream PEBCOM export that added case 791897762
# ("Aristobulo del Valle 1707").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing row 7 (and everything below it) down by one row.
$ws.Rows(7).Insert()

# The source feed stores Caso / Fecha / Direccion / Comuna / OT / Proveedor /
# Estado / Observaciones / Tipo de tarea / Equipo / Tipo de Elemento /
# Operacion / Zona as plain text, even when the text looks numeric or like a
# date. Force the text number format on the new row before writing the
# values so Excel doesn't auto-convert them to numbers/dates.
$ws.Range("A7:H7").NumberFormat = "@"
$ws.Range("J7:L7").NumberFormat = "@"
$ws.Range("O7:P7").NumberFormat = "@"

$ws.Range("A7").Value = "791897762"
$ws.Range("B7").Value = "8/2/2024"
$ws.Range("C7").Value = "Aristobulo del Valle 1707"
$ws.Range("D7").Value = "4"
$ws.Range("E7").Value = "791897762"
$ws.Range("F7").Value = "PEBCOM"
$ws.Range("G7").Value = "Pendiente"
$ws.Range("H7").Value = "Pendiente"
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = "Cambio"
$ws.Range("K7").Value = "Sin equipos"
$ws.Range("M7").Value = -58.375312
$ws.Range("N7").Value = -34.636076
$ws.Range("O7").Value = "San Telmo"
$ws.Range("P7").Value = "Capital Sur"
